$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30750
$ws.Range("J3").Value = 30750
$ws.Range("N3").Value = -30978
$ws.Range("L3").Value = 30750
$ws.Range("N4").Value = -1166.3333
$ws.Range("H4").Value = 528.375
$ws.Range("M4").Value = -168.4
$ws.Range("I4").Value = 282.4
$ws.Range("K4").Value = 282.4
$ws.Range("J4").Value = 938.3333
$ws.Range("L4").Value = 938.3333
$ws.Range("H15").Value = 2691.0657
$ws.Range("I15").Value = 2691.0657
$ws.Range("K15").Value = 8073.1971
$ws.Range("M15").Value = -7904.1971
$ws.Range("J17").Value = 2378.4285
$ws.Range("L17").Value = 7135.2855
$ws.Range("H17").Value = 2252.7646
$ws.Range("N17").Value = -7471.2855
$ws.Range("K28").Value = 453.33334
$ws.Range("I28").Value = 453.33334
$ws.Range("H28").Value = 523.2857
$ws.Range("M28").Value = 31.66665999999998
$ws.Range("I41").Value = 243.58333
$ws.Range("J41").Value = 507
$ws.Range("K41").Value = 243.58333
$ws.Range("M41").Value = 196.41667
$ws.Range("L41").Value = 507
$ws.Range("H41").Value = 389.92593
$ws.Range("N41").Value = -1387
$ws.Range("M53").Value = 396.33333
$ws.Range("I53").Value = 240.66667
$ws.Range("K53").Value = 240.66667
$ws.Range("H53").Value = 319.36365
$ws.Range("H62").Value = 2998.2
$ws.Range("M62").Value = -2375
$ws.Range("K62").Value = 2999
$ws.Range("I62").Value = 2999
$ws.Range("M64").Value = -3346
$ws.Range("K64").Value = 3594
$ws.Range("H64").Value = 6864.6665
$ws.Range("I64").Value = 3594
$ws.Range("I65").Value = 2999
$ws.Range("H65").Value = 2998.2
$ws.Range("K65").Value = 14995
$ws.Range("M65").Value = -11875
$ws.Range("H67").Value = 6864.6665
$ws.Range("I67").Value = 3594
$ws.Range("K67").Value = 3594
$ws.Range("M67").Value = -2736
$ws.Range("K86").Value = 181824740
$ws.Range("H86").Value = 111117720
$ws.Range("I86").Value = 181824740
$ws.Range("M86").Value = -181823617
$ws.Range("I89").Value = 181824740
$ws.Range("K89").Value = 909123700
$ws.Range("M89").Value = -909118084
$ws.Range("H89").Value = 111117720
$ws.Range("J98").Value = 1565
$ws.Range("N98").Value = -4561
$ws.Range("L98").Value = 1565
$ws.Range("H98").Value = 1398.7715
$ws.Range("J102").Value = 30750
$ws.Range("H102").Value = 30750
$ws.Range("N102").Value = -37240
$ws.Range("L102").Value = 30750
$ws.Range("I106").Value = 3959
$ws.Range("M106").Value = -3328
$ws.Range("H106").Value = 4252.7617
$ws.Range("K106").Value = 3959
$ws.Range("J116").Value = 5437.7
$ws.Range("L116").Value = 5437.7
$ws.Range("K116").Value = 6255.048
$ws.Range("N116").Value = -12321.7
$ws.Range("M116").Value = -2813.048
$ws.Range("I116").Value = 6255.048
$ws.Range("H116").Value = 5991.387
$ws.Range("L121").Value = 8464.799999999999
$ws.Range("N121").Value = -11958.8
$ws.Range("H121").Value = 3693.8333
$ws.Range("M121").Value = -22418
$ws.Range("I121").Value = 8055
$ws.Range("J121").Value = 2821.6
$ws.Range("K121").Value = 24165
$ws.Range("N122").Value = -9595
$ws.Range("J122").Value = 1565
$ws.Range("L122").Value = 4695
$ws.Range("H122").Value = 1398.7715
$ws.Range("J125").Value = 7896.625
$ws.Range("N125").Value = -75989.625
$ws.Range("L125").Value = 71069.625
$ws.Range("H125").Value = 4428.8423
$ws.Range("M132").Value = -10704.2438
$ws.Range("I132").Value = 4411.4146
$ws.Range("H132").Value = 6562.778
$ws.Range("K132").Value = 13234.2438
$ws.Range("I135").Value = 414.5
$ws.Range("L135").Value = 33750
$ws.Range("H135").Value = 1248.375
$ws.Range("K135").Value = 3730.5
$ws.Range("N135").Value = -38820
$ws.Range("M135").Value = -1195.5
$ws.Range("J135").Value = 3750
$ws.Range("L136").Value = 110253.27
$ws.Range("J136").Value = 110253.27
$ws.Range("H136").Value = 110253.27
$ws.Range("N136").Value = -120453.27
$ws.Range("I137").Value = 3858.925
$ws.Range("J137").Value = 1930.9
$ws.Range("N137").Value = -10892.7
$ws.Range("H137").Value = 3473.32
$ws.Range("K137").Value = 11576.775
$ws.Range("M137").Value = -9026.775000000001
$ws.Range("L137").Value = 5792.700000000001
$ws.Range("L138").Value = 12958.2
$ws.Range("J138").Value = 4319.4
$ws.Range("N138").Value = -23238.2
$ws.Range("L140").Value = 75160.75
$ws.Range("N140").Value = -85520.75
$ws.Range("H140").Value = 75160.75
$ws.Range("J140").Value = 75160.75
$ws.Range("I141").Value = 18755.846
$ws.Range("H141").Value = 12620.714
$ws.Range("K141").Value = 56267.538
$ws.Range("M141").Value = -51087.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5003000
$ws.Range("K28").Value = 5392.3335
$ws.Range("I28").Value = 5392.3335
$ws.Range("H28").Value = 20772.125
$ws.Range("M28").Value = -5200.3335
$ws.Range("H32").Value = 1050.4
$ws.Range("I32").Value = 1050.4
$ws.Range("K32").Value = 1050.4
$ws.Range("M32").Value = -763.4000000000001
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30972
$ws.Range("M45").Value = -13369
$ws.Range("H45").Value = 13328.777
$ws.Range("I45").Value = 13746
$ws.Range("K45").Value = 13746
$ws.Range("K61").Value = 5451.522
$ws.Range("I61").Value = 5451.522
$ws.Range("H61").Value = 6057.7085
$ws.Range("M61").Value = -5239.522
$ws.Range("K74").Value = 1556.9688
$ws.Range("M74").Value = -682.9688000000001
$ws.Range("I74").Value = 1556.9688
$ws.Range("L74").Value = 3632.2856
$ws.Range("J74").Value = 3632.2856
$ws.Range("N74").Value = -5380.2856
$ws.Range("H74").Value = 1929.4615
$ws.Range("K77").Value = 7784.844000000001
$ws.Range("M77").Value = -3416.844000000001
$ws.Range("N77").Value = -26897.428
$ws.Range("L77").Value = 18161.428
$ws.Range("I77").Value = 1556.9688
$ws.Range("H77").Value = 1929.4615
$ws.Range("J77").Value = 3632.2856
$ws.Range("I99").Value = 5392.3335
$ws.Range("H99").Value = 20772.125
$ws.Range("M99").Value = -2397.3335
$ws.Range("K99").Value = 5392.3335
$ws.Range("J118").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("L118").Value = 0
$ws.Range("M122").Value = -8005.856800000001
$ws.Range("I122").Value = 3485.2856
$ws.Range("K122").Value = 10455.8568
$ws.Range("H122").Value = 3674.5
$ws.Range("H130").Value = 77564.5
$ws.Range("J130").Value = 77564.5
$ws.Range("N130").Value = -87604.5
$ws.Range("L130").Value = 77564.5
$ws.Range("J132").Value = 3496.4614
$ws.Range("N132").Value = -15549.3842
$ws.Range("M132").Value = -5657.6666
$ws.Range("L132").Value = 10489.3842
$ws.Range("I132").Value = 2729.2222
$ws.Range("H132").Value = 3050.9678
$ws.Range("K132").Value = 8187.6666
$ws.Range("M136").Value = -13804.566
$ws.Range("I136").Value = 5451.522
$ws.Range("K136").Value = 16354.566
$ws.Range("H136").Value = 6057.7085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J7").Value = 2997.5
$ws.Range("L7").Value = 2997.5
$ws.Range("N7").Value = -3223.5
$ws.Range("H7").Value = 2498.75
$ws.Range("K86").Value = 3141.75
$ws.Range("H86").Value = 3358.85
$ws.Range("I86").Value = 3141.75
$ws.Range("M86").Value = -2018.75
$ws.Range("I89").Value = 3141.75
$ws.Range("K89").Value = 15708.75
$ws.Range("M89").Value = -10092.75
$ws.Range("H89").Value = 3358.85
$ws.Range("I94").Value = 946.3077
$ws.Range("K94").Value = 946.3077
$ws.Range("H94").Value = 1966.7142
$ws.Range("M94").Value = -495.3077
$ws.Range("I96").Value = 10509.333
$ws.Range("H96").Value = 10509.333
$ws.Range("K96").Value = 10509.333
$ws.Range("M96").Value = -7763.333000000001
$ws.Range("N99").Value = -8013.6
$ws.Range("I99").Value = 2941.182
$ws.Range("H99").Value = 3590.0625
$ws.Range("J99").Value = 5017.6
$ws.Range("L99").Value = 5017.6
$ws.Range("M99").Value = -1443.182
$ws.Range("K99").Value = 2941.182
$ws.Range("H105").Value = 2435.238
$ws.Range("I105").Value = 2375.7896
$ws.Range("K105").Value = 2375.7896
$ws.Range("M105").Value = -628.7896000000001
$ws.Range("I107").Value = 3062.1191
$ws.Range("L107").Value = 3801.6875
$ws.Range("H107").Value = 3266.138
$ws.Range("N107").Value = -7641.6875
$ws.Range("K107").Value = 3062.1191
$ws.Range("M107").Value = -1142.1191
$ws.Range("J107").Value = 3801.6875
$ws.Range("M134").Value = -4879.1535
$ws.Range("K134").Value = 7414.1535
$ws.Range("I134").Value = 2471.3845
$ws.Range("H134").Value = 2838.634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I2").Value = 15.666667
$ws.Range("H2").Value = 48.2
$ws.Range("L2").Value = 97
$ws.Range("J2").Value = 97
$ws.Range("M2").Value = 97.333333
$ws.Range("K2").Value = 15.666667
$ws.Range("N2").Value = -323
$ws.Range("H4").Value = 10000
$ws.Range("M4").Value = -9888
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("L5").Value = 1298.375
$ws.Range("M5").Value = -2359
$ws.Range("I5").Value = 2471
$ws.Range("H5").Value = 1949.8334
$ws.Range("N5").Value = -1522.375
$ws.Range("K5").Value = 2471
$ws.Range("J5").Value = 1298.375
$ws.Range("L6").Value = 998
$ws.Range("M6").Value = -500886.3
$ws.Range("N6").Value = -1224
$ws.Range("K6").Value = 500999.3
$ws.Range("I6").Value = 500999.3
$ws.Range("H6").Value = 455544.62
$ws.Range("J6").Value = 998
$ws.Range("N10").Value = -5276.75
$ws.Range("J10").Value = 4998.75
$ws.Range("M10").Value = -306.63635
$ws.Range("K10").Value = 445.63635
$ws.Range("L10").Value = 4998.75
$ws.Range("I10").Value = 445.63635
$ws.Range("H10").Value = 1659.8
$ws.Range("L11").Value = 799.3333
$ws.Range("J11").Value = 799.3333
$ws.Range("N11").Value = -1079.3333
$ws.Range("H11").Value = 799.3333
$ws.Range("M12").Value = -5431.6665
$ws.Range("K12").Value = 5601.6665
$ws.Range("N12").Value = -10000340
$ws.Range("J12").Value = 10000000
$ws.Range("L12").Value = 10000000
$ws.Range("I12").Value = 5601.6665
$ws.Range("H12").Value = 2504201.2
$ws.Range("K31").Value = 2108.875
$ws.Range("N31").Value = -5681.6772
$ws.Range("H31").Value = 3576.6033
$ws.Range("J31").Value = 5091.6772
$ws.Range("M31").Value = -1813.875
$ws.Range("L31").Value = 5091.6772
$ws.Range("I31").Value = 2108.875
$ws.Range("J34").Value = 5091.6772
$ws.Range("K34").Value = 2108.875
$ws.Range("L34").Value = 5091.6772
$ws.Range("M34").Value = -1906.875
$ws.Range("N34").Value = -5495.6772
$ws.Range("I34").Value = 2108.875
$ws.Range("H34").Value = 3576.6033
$ws.Range("H58").Value = 2806.647
$ws.Range("K58").Value = 1844.8182
$ws.Range("L58").Value = 4570
$ws.Range("N58").Value = -4976
$ws.Range("I58").Value = 1844.8182
$ws.Range("J58").Value = 4570
$ws.Range("M58").Value = -1641.8182
$ws.Range("H69").Value = 19999
$ws.Range("K69").Value = 19999
$ws.Range("I69").Value = 19999
$ws.Range("M69").Value = -19250
$ws.Range("I72").Value = 19999
$ws.Range("K72").Value = 59997
$ws.Range("H72").Value = 19999
$ws.Range("M72").Value = -56253
$ws.Range("L92").Value = 38499
$ws.Range("H92").Value = 38499
$ws.Range("J92").Value = 38499
$ws.Range("N92").Value = -43491
$ws.Range("J132").Value = 6348.3335
$ws.Range("N132").Value = -24105.0005
$ws.Range("M132").Value = -12198.971
$ws.Range("L132").Value = 19045.0005
$ws.Range("I132").Value = 4909.657
$ws.Range("H132").Value = 5120.1953
$ws.Range("K132").Value = 14728.971
$ws.Range("M134").Value = -5015.8125
$ws.Range("K134").Value = 7550.8125
$ws.Range("I134").Value = 2516.9375
$ws.Range("H134").Value = 2753.17
$ws.Range("M136").Value = -2984.4546
$ws.Range("L136").Value = 13710
$ws.Range("I136").Value = 1844.8182
$ws.Range("K136").Value = 5534.4546
$ws.Range("J136").Value = 4570
$ws.Range("H136").Value = 2806.647
$ws.Range("N136").Value = -18810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10833.167
$ws.Range("J39").Value = 11772.546
$ws.Range("L39").Value = 35317.638
$ws.Range("N39").Value = -35905.638
$ws.Range("L107").Value = 3046.66662
$ws.Range("H107").Value = 1333
$ws.Range("N107").Value = -6886.66662
$ws.Range("J107").Value = 1015.55554
$ws.Range("J117").Value = 2199.7778
$ws.Range("N117").Value = -13483.3334
$ws.Range("H117").Value = 3635.2727
$ws.Range("L117").Value = 6599.3334
$ws.Range("N129").Value = -14662
$ws.Range("J129").Value = 1554
$ws.Range("H129").Value = 1197.7368
$ws.Range("L129").Value = 4662
$ws.Range("K131").Value = 3678175.2
$ws.Range("I131").Value = 1226058.4
$ws.Range("L131").Value = 5207.1
$ws.Range("J131").Value = 1735.7
$ws.Range("H131").Value = 865963.5
$ws.Range("N131").Value = -15287.1
$ws.Range("M131").Value = -3673135.2
$ws.Range("J132").Value = 3263.6667
$ws.Range("N132").Value = -34433.0003
$ws.Range("M132").Value = -6890.75
$ws.Range("L132").Value = 29373.0003
$ws.Range("I132").Value = 1046.75
$ws.Range("H132").Value = 2581.5386
$ws.Range("K132").Value = 9420.75
$ws.Range("J137").Value = 6513.5
$ws.Range("N137").Value = -29740.5
$ws.Range("H137").Value = 4014.375
$ws.Range("L137").Value = 19540.5
$ws.Range("I141").Value = 5985.8
$ws.Range("H141").Value = 5985.8
$ws.Range("K141").Value = 17957.4
$ws.Range("M141").Value = -12777.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K9").Value = 407
$ws.Range("I9").Value = 407
$ws.Range("M9").Value = -237
$ws.Range("H9").Value = 407
$ws.Range("M11").Value = -2824861.8
$ws.Range("I11").Value = 2825000.8
$ws.Range("L11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H11").Value = 2825000.8
$ws.Range("K11").Value = 2825000.8
$ws.Range("M13").Value = -88.25
$ws.Range("K13").Value = 227.25
$ws.Range("I13").Value = 227.25
$ws.Range("H13").Value = 302
$ws.Range("J17").Value = 999
$ws.Range("L17").Value = 999
$ws.Range("H17").Value = 1237.4615
$ws.Range("N17").Value = -1335
$ws.Range("J101").Value = 47922.5
$ws.Range("H101").Value = 47922.5
$ws.Range("L101").Value = 47922.5
$ws.Range("N101").Value = -54412.5
$ws.Range("M102").Value = -2571.265
$ws.Range("K102").Value = 4193.265
$ws.Range("I102").Value = 4193.265
$ws.Range("H102").Value = 5163.302
$ws.Range("I107").Value = 864.8333
$ws.Range("L107").Value = 1182.125
$ws.Range("H107").Value = 1046.1428
$ws.Range("N107").Value = -5022.125
$ws.Range("K107").Value = 864.8333
$ws.Range("M107").Value = 1055.1667
$ws.Range("J107").Value = 1182.125
$ws.Range("K113").Value = 2533.75
$ws.Range("N113").Value = -8854.666499999999
$ws.Range("J113").Value = 4514.6665
$ws.Range("M113").Value = -363.75
$ws.Range("L113").Value = 4514.6665
$ws.Range("H113").Value = 3722.3
$ws.Range("I113").Value = 2533.75
$ws.Range("J117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M122").Value = -5353.599999999999
$ws.Range("I122").Value = 2601.2
$ws.Range("K122").Value = 7803.599999999999
$ws.Range("H122").Value = 3027.65
$ws.Range("H126").Value = 5674.857
$ws.Range("M126").Value = -10269.0005
$ws.Range("I126").Value = 4246.3335
$ws.Range("K126").Value = 12739.0005
$ws.Range("J132").Value = 5243.625
$ws.Range("N132").Value = -20790.875
$ws.Range("M132").Value = -7832.960000000001
$ws.Range("L132").Value = 15730.875
$ws.Range("I132").Value = 3454.32
$ws.Range("H132").Value = 3888.0908
$ws.Range("K132").Value = 10362.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M11").Value = -5758
$ws.Range("I11").Value = 5898
$ws.Range("L11").Value = 5800
$ws.Range("J11").Value = 5800
$ws.Range("N11").Value = -6080
$ws.Range("H11").Value = 5873.5
$ws.Range("K11").Value = 5898
$ws.Range("H35").Value = 6662.8335
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -6326.8335
$ws.Range("K35").Value = 6662.8335
$ws.Range("J35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("I35").Value = 6662.8335
$ws.Range("K61").Value = 3196.484
$ws.Range("I61").Value = 3196.484
$ws.Range("H61").Value = 3148.1562
$ws.Range("M61").Value = -2994.484
$ws.Range("K113").Value = 3196.484
$ws.Range("M113").Value = -1026.484
$ws.Range("H113").Value = 3148.1562
$ws.Range("I113").Value = 3196.484
$ws.Range("L115").Value = 175000
$ws.Range("J115").Value = 175000
$ws.Range("N115").Value = -177350
$ws.Range("H115").Value = 175000
$ws.Range("N122").ClearContents()
$ws.Range("J122").Value = 0
$ws.Range("M122").Value = -59299
$ws.Range("L122").Value = 0
$ws.Range("I122").Value = 20583
$ws.Range("K122").Value = 61749
$ws.Range("H122").Value = 20583
$ws.Range("J132").Value = 10523.958
$ws.Range("N132").Value = -36631.874
$ws.Range("M132").Value = -149908.448
$ws.Range("L132").Value = 31571.874
$ws.Range("I132").Value = 50812.816
$ws.Range("H132").Value = 29792.543
$ws.Range("K132").Value = 152438.448
$ws.Range("M136").Value = -33780537
$ws.Range("I136").Value = 11261029
$ws.Range("K136").Value = 33783087
$ws.Range("H136").Value = 6676792.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N4").Value = -22692
$ws.Range("H4").Value = 15477.333
$ws.Range("M4").Value = -1387
$ws.Range("I4").Value = 1500
$ws.Range("K4").Value = 1500
$ws.Range("J4").Value = 22466
$ws.Range("L4").Value = 22466
$ws.Range("L5").Value = 6250.5
$ws.Range("M5").ClearContents()
$ws.Range("I5").Value = 0
$ws.Range("H5").Value = 6250.5
$ws.Range("N5").Value = -6474.5
$ws.Range("K5").Value = 0
$ws.Range("J5").Value = 6250.5
$ws.Range("I7").Value = 50000
$ws.Range("K7").Value = 50000
$ws.Range("M7").Value = -49887
$ws.Range("H7").Value = 50000
$ws.Range("K9").Value = 14664.667
$ws.Range("I9").Value = 14664.667
$ws.Range("M9").Value = -14524.667
$ws.Range("H9").Value = 15998.5
$ws.Range("J29").Value = 1405.5
$ws.Range("I29").Value = 4499.75
$ws.Range("H29").Value = 3468.3333
$ws.Range("L29").Value = 1405.5
$ws.Range("M29").Value = -4209.75
$ws.Range("N29").Value = -1985.5
$ws.Range("K29").Value = 4499.75
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("N103").Value = -36019.668
$ws.Range("J103").Value = 33675.668
$ws.Range("H103").Value = 33675.668
$ws.Range("L103").Value = 33675.668
$ws.Range("N113").Value = -7338.5
$ws.Range("J113").Value = 999.5
$ws.Range("L113").Value = 2998.5
$ws.Range("H113").Value = 653.8570999999999
$ws.Range("N122").Value = -96697.75
$ws.Range("J122").Value = 30599.25
$ws.Range("M122").Value = -14472.5005
$ws.Range("L122").Value = 91797.75
$ws.Range("I122").Value = 5640.8335
$ws.Range("K122").Value = 16922.5005
$ws.Range("H122").Value = 9206.321
$ws.Range("M127").ClearContents()
$ws.Range("L127").Value = 66666.664
$ws.Range("J127").Value = 66666.664
$ws.Range("I127").Value = 0
$ws.Range("N127").Value = -76586.664
$ws.Range("K127").Value = 0
$ws.Range("H127").Value = 66666.664
$ws.Range("J132").Value = 7602.8184
$ws.Range("N132").Value = -27868.4552
$ws.Range("M132").Value = -5054
$ws.Range("L132").Value = 22808.4552
$ws.Range("I132").Value = 2528
$ws.Range("H132").Value = 3906.3457
$ws.Range("K132").Value = 7584
$ws.Range("M136").Value = -2971.683
$ws.Range("L136").Value = 16015.8
$ws.Range("I136").Value = 1840.561
$ws.Range("K136").Value = 5521.683
$ws.Range("J136").Value = 5338.6
$ws.Range("H136").Value = 2526.451
$ws.Range("N136").Value = -21115.8
$ws.Range("J141").Value = 158263.5
$ws.Range("N141").Value = -168623.5
$ws.Range("L141").Value = 158263.5
$ws.Range("H141").Value = 158263.5
